# LC Case Study - final changes
# Slide 2, Shape "TextBox 4": grow the textbox and append two new
# paragraphs (a blank spacer + the "Problem to be solved" note in a
# muted gray Consolas run) after the existing objective paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(1)

# --- Resize the shape (EMU -> points, PowerPoint COM uses points) ---
$EMU_PER_POINT = 12700
$shape.Height = 2554545 / $EMU_PER_POINT

# --- Append the new paragraphs to the text body ---
$tr = $shape.TextFrame.TextRange

$problemText = "Problem to be solved: Identify patterns based on the columns which would be used to reject/deny loans, reduce the loan amount or charge higher interest rate."

$beforeLen = $tr.Length
# One blank paragraph, then the new "Problem to be solved" paragraph,
# then a trailing blank paragraph (the original trailing paragraph mark
# is preserved after this insertion).
$insertText = "`r`r" + $problemText + "`r`r"
$tr.InsertAfter($insertText) | Out-Null

# Grab just the inserted "Problem to be solved" run so we can restyle it
# without touching anything else.
$problemStart = $beforeLen + 3
$problemRange = $tr.Characters($problemStart, $problemText.Length)

$problemRange.Font.Name = "Consolas"
$problemRange.Font.Bold = 0
$problemRange.Font.Color.RGB = 0xCCCCCC
# Forces an explicit (empty) effect list, matching the authored markup.
$problemRange.Font.Shadow = 0
